# ---------------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# The workbook used to have a single quarterly sheet named "2022-Q3". This
# change introduces a new "2022-Q4" sheet (with refreshed numbers) while
# keeping a "2022-Q3" sheet around with the data it already had.
#
# Concretely or equivalently:
#   1. The existing "2022-Q3" sheet is duplicated, so the duplicate keeps the
#      original Q3 numbers/formatting untouched and becomes the new
#      "2022-Q3" sheet.
#   2. The original "2022-Q3" sheet is renamed to "2022-Q4" and its holdings
#      figures are refreshed with the new quarter's numbers.
#   3. The summary ("总计") sheet gets a new row for the 2022-Q3 entry (the
#      2022-Q4 entry reuses the row that used to say "2022-Q3").
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Helper: write $text into $sheet's $addr cell as a genuine text value (even
# when $text looks like a number), without leaving any extra formatting on
# the destination cell. We stage the text on a scratch cell that is forced
# to text format, copy just its value onto the destination, then discard the
# scratch column again.
function Set-TextValue($sheet, $addr, $text) {
    $scratch = $sheet.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $sheet.Range($addr).PasteSpecial(-4163)   # xlPasteValues
}

# --- Step 1: duplicate "2022-Q3" so the old data survives on its own sheet ---
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $q3)

# --- Step 2: turn the original sheet into "2022-Q4" with refreshed figures ---
# (Rename the original out of the way first so the duplicate can reclaim the
# "2022-Q3" name.)
$q4 = $q3
$q4.Name = "2022-Q4"

$q3Old = $wb.Worksheets.Item("2022-Q3 (2)")
$q3Old.Name = "2022-Q3"

# Match the header/first-column styling used on the summary sheet.
$summary = $wb.Worksheets.Item("总计")
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)        # xlPasteFormats
$summary.Range("A2").Copy()
$q4.Range("A2:A3").PasteSpecial(-4122)        # xlPasteFormats

Set-TextValue $q4 "D2" "0.24"
Set-TextValue $q4 "E2" "94.45"
Set-TextValue $q4 "F2" "1.46"
Set-TextValue $q4 "G2" "0.0035"
$q4.Range("H2").Value = 6

Set-TextValue $q4 "D3" "0.05"
Set-TextValue $q4 "E3" "94.45"
Set-TextValue $q4 "F3" "1.46"
Set-TextValue $q4 "G3" "0.0007"
$q4.Range("H3").Value = 6

$q4.Range("ZZ1").EntireColumn.Delete()

# --- Step 3: update the summary ("总计") sheet ---
$summary.Range("B2").Value = "2022-Q4"

$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)      # xlPasteFormats
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0
